$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$lastSheet = $wb.Worksheets.Item("2021-Q3")

# --- 1. Create the new "2022-Q4" worksheet right after "总计" by
#        duplicating the "2022-Q3" sheet (keeps identical page setup /
#        sheetPr / styles) and then replacing its contents. ---
$q3Sheet.Copy($null, $summarySheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Drop all the rows that came along with the copied "2022-Q3" data,
# keeping only the header row.
$newSheet.Range("A4:H50").EntireRow.Delete()

# --- 2. Fill in the new "2022-Q4" sheet contents ---
# Columns B-G hold text values in the source data (fund codes, names and
# numbers-as-text), so force text formatting before writing them in so
# leading zeros / decimal text are preserved verbatim.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "011174"
$newSheet.Range("C2").Value = "中庚价值品质一年持有期混合"
$newSheet.Range("D2").Value = "67.05"
$newSheet.Range("E2").Value = "93.59"
$newSheet.Range("F2").Value = "2.51"
$newSheet.Range("G2").Value = "1.6830"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "004099"
$newSheet.Range("C3").Value = "前海开源沪港深景气行业精选灵活配置混合"
$newSheet.Range("D3").Value = "0.39"
$newSheet.Range("E3").Value = "88.59"
$newSheet.Range("F3").Value = "4.97"
$newSheet.Range("G3").Value = "0.0194"
$newSheet.Range("H3").Value = 7

# Drop the text number-format again so the cells end up with the same
# "no explicit style" look as the rest of the workbook's data cells.
$newSheet.Range("B2:G3").Style = "Normal"

# --- 3. Update the "总计" sheet: shift the quarter rows down by one and
#        insert the new 2022-Q4 figures at the top ---
# Row 7 is a brand new row, so first copy the "A" column style down from
# row 6 before filling it in.
$summarySheet.Range("A6").Copy()
$summarySheet.Range("A7").PasteSpecial(-4122)

$summarySheet.Range("A7").Value = 5
$summarySheet.Range("B7").Value = "2021-Q3"
$summarySheet.Range("C7").Value = 2
$summarySheet.Range("D7").Value = 1.36

$summarySheet.Range("A6").Value = 4
$summarySheet.Range("B6").Value = "2021-Q4"
$summarySheet.Range("C6").Value = 2
$summarySheet.Range("D6").Value = 0.28

$summarySheet.Range("A5").Value = 3
$summarySheet.Range("B5").Value = "2022-Q1"
$summarySheet.Range("C5").Value = 6
$summarySheet.Range("D5").Value = 1.25

$summarySheet.Range("A4").Value = 2
$summarySheet.Range("B4").Value = "2022-Q2"
$summarySheet.Range("C4").Value = 14
$summarySheet.Range("D4").Value = 1.4

$summarySheet.Range("A3").Value = 1
$summarySheet.Range("B3").Value = "2022-Q3"
$summarySheet.Range("C3").Value = 49
$summarySheet.Range("D3").Value = 8.800000000000001

$summarySheet.Range("A2").Value = 0
$summarySheet.Range("B2").Value = "2022-Q4"
$summarySheet.Range("C2").Value = 2
$summarySheet.Range("D2").Value = 1.7

# --- 4. Restore the original tab-selection state: keep "2021-Q3" marked
#        as the selected tab and "总计" as the active workbook tab. ---
$lastSheet.Select($false)
$summarySheet.Activate()
